# Regenerate save_data G column ("K" = strikeouts) with recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..36 (column G), replacing the old Strike# derived figures.
$kValues = @{
    2  = 3
    3  = 4
    4  = 10
    5  = 8
    6  = 6
    7  = 7
    8  = 4
    9  = 11
    10 = 4
    11 = 3
    12 = 8
    13 = 5
    14 = 3
    15 = 8
    16 = 0
    17 = 5
    18 = 9
    19 = 7
    20 = 10
    21 = 5
    22 = 6
    23 = 9
    24 = 5
    25 = 0
    26 = 8
    27 = 8
    28 = 4
    29 = 17
    30 = 5
    31 = 4
    32 = 9
    33 = 5
    34 = 1
    35 = 4
    36 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
